$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.071944773141773
$ws.Range("D2").Value = 1.065553645872923
$ws.Range("E2").Value = 1.08564895773982
$ws.Range("F2").Value = 1.09286200330311
$ws.Range("I2").Value = 1.05260349743665
$ws.Range("J2").Value = 1.076866809375162
$ws.Range("K2").Value = 1.068266574889631
$ws.Range("L2").Value = 1.088308842025091
$ws.Range("M2").Value = 1.095503360871276
$ws.Range("N2").Value = 1.078396083762825

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.073394417564274
$ws.Range("D3").Value = 1.066647634065784
$ws.Range("E3").Value = 1.087126625846894
$ws.Range("F3").Value = 1.09448179813807
$ws.Range("I3").Value = 1.05306110335575
$ws.Range("J3").Value = 1.077972225130458
$ws.Range("K3").Value = 1.069175288265812
$ws.Range("L3").Value = 1.089604159099629
$ws.Range("M3").Value = 1.096941821713494
$ws.Range("N3").Value = 1.079503069335286

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.074331094205834
$ws.Range("D4").Value = 1.067354243437833
$ws.Range("E4").Value = 1.08808179810815
$ws.Range("F4").Value = 1.095529171945591
$ws.Range("I4").Value = 1.053355264646675
$ws.Range("J4").Value = 1.078685694945682
$ws.Range("K4").Value = 1.069761419457152
$ws.Range("L4").Value = 1.090440816118226
$ws.Range("M4").Value = 1.097871363299531
$ws.Range("N4").Value = 1.080217552359484

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.074724558930572
$ws.Range("D5").Value = 1.067651001250119
$ws.Range("E5").Value = 1.088483124836999
$ws.Range("F5").Value = 1.095969317460608
$ws.Range("I5").Value = 1.053478467872524
$ws.Range("J5").Value = 1.078985209890603
$ws.Range("K5").Value = 1.070007386106912
$ws.Range("L5").Value = 1.090792193550103
$ws.Range("M5").Value = 1.098261852506322
$ws.Range("N5").Value = 1.080517492649981

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.074790605162818
$ws.Range("D6").Value = 1.067700810678897
$ws.Range("E6").Value = 1.088550496194291
$ws.Range("F6").Value = 1.096043210042059
$ws.Range("I6").Value = 1.053499127201512
$ws.Range("J6").Value = 1.079035474831022
$ws.Range("K6").Value = 1.07004865913612
$ws.Range("L6").Value = 1.090851170826319
$ws.Range("M6").Value = 1.098327400540125
$ws.Range("N6").Value = 1.080567828972381

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.074336352926749
$ws.Range("D7").Value = 1.06735820990678
$ws.Range("E7").Value = 1.088087161542272
$ws.Range("F7").Value = 1.095535053852327
$ws.Range("I7").Value = 1.053356912707372
$ws.Range("J7").Value = 1.078689698753554
$ws.Range("K7").Value = 1.069764707810813
$ws.Range("L7").Value = 1.090445512621158
$ws.Range("M7").Value = 1.097876582167184
$ws.Range("N7").Value = 1.080221561853222

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.072434968819754
$ws.Range("D8").Value = 1.06592363069829
$ws.Range("E8").Value = 1.086148548919176
$ws.Range("F8").Value = 1.093409578449149
$ws.Range("I8").Value = 1.052758550671624
$ws.Range("J8").Value = 1.077240767258536
$ws.Range("K8").Value = 1.068574067899686
$ws.Range("L8").Value = 1.088746914564413
$ws.Range("M8").Value = 1.095989755856039
$ws.Range("N8").Value = 1.078770572709287

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.069073920698476
$ws.Range("D9").Value = 1.063385768519854
$ws.Range("E9").Value = 1.08272468937009
$ws.Range("F9").Value = 1.089658217309772
$ws.Range("I9").Value = 1.051689194345187
$ws.Range("J9").Value = 1.07467348790553
$ws.Range("K9").Value = 1.066461524273191
$ws.Range("L9").Value = 1.08574201807699
$ws.Range("M9").Value = 1.092655147966015
$ws.Range("N9").Value = 1.076199647525128

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.066825693902878
$ws.Range("D10").Value = 1.061686899753212
$ws.Range("E10").Value = 1.080436483882926
$ws.Range("F10").Value = 1.087152823993632
$ws.Range("I10").Value = 1.050966079381606
$ws.Range("J10").Value = 1.072952181588993
$ws.Range("K10").Value = 1.065043164462324
$ws.Range("L10").Value = 1.083730483649676
$ws.Range("M10").Value = 1.090425091971111
$ws.Range("N10").Value = 1.074475896756185

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.06585030648002
$ws.Range("D11").Value = 1.060949558945749
$ws.Range("E11").Value = 1.079444239206645
$ws.Range("F11").Value = 1.086066796232749
$ws.Range("I11").Value = 1.050650507337453
$ws.Range("J11").Value = 1.072204443389148
$ws.Range("K11").Value = 1.064426569507301
$ws.Range("L11").Value = 1.082857422504758
$ws.Range("M11").Value = 1.089457704933295
$ws.Range("N11").Value = 1.07372709668233

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.06548771190319
$ws.Range("D12").Value = 1.06067541408102
$ws.Range("E12").Value = 1.079075451645404
$ws.Range("F12").Value = 1.085663211739821
$ws.Range("I12").Value = 1.050532917524617
$ws.Range("J12").Value = 1.071926333063435
$ws.Range("K12").Value = 1.064197167634204
$ws.Range("L12").Value = 1.08253281298171
$ws.Range("M12").Value = 1.089098102208616
$ws.Range("N12").Value = 1.073448591408054

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.065565503029176
$ws.Range("D13").Value = 1.060734231102331
$ws.Range("E13").Value = 1.07915456807951
$ws.Range("F13").Value = 1.085749790515791
$ws.Range("I13").Value = 1.050558157833711
$ws.Range("J13").Value = 1.071986005384659
$ws.Range("K13").Value = 1.064246391987952
$ws.Range("L13").Value = 1.082602457230898
$ws.Range("M13").Value = 1.089175250655764
$ws.Range("N13").Value = 1.073508348470819

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.065820340309101
$ws.Range("D14").Value = 1.060926903455277
$ws.Range("E14").Value = 1.079413759721132
$ws.Range("F14").Value = 1.086033439646195
$ws.Range("I14").Value = 1.050640794938851
$ws.Range("J14").Value = 1.072181462227554
$ws.Range("K14").Value = 1.064407614675921
$ws.Range("L14").Value = 1.082830596663343
$ws.Range("M14").Value = 1.089427985657098
$ws.Range("N14").Value = 1.073704082884849

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.065977314993051
$ws.Range("D15").Value = 1.061045580153263
$ws.Range("E15").Value = 1.079573426365808
$ws.Range("F15").Value = 1.086208180403114
$ws.Range("I15").Value = 1.050691660964264
$ws.Range("J15").Value = 1.072301840810072
$ws.Range("K15").Value = 1.064506899973616
$ws.Range("L15").Value = 1.082971118841321
$ws.Range("M15").Value = 1.089583667736538
$ws.Range("N15").Value = 1.073824632418762

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.066890386622664
$ws.Range("D16").Value = 1.061735797953572
$ws.Range("E16").Value = 1.080502304955713
$ws.Range("F16").Value = 1.087224874557873
$ws.Range("I16").Value = 1.050986970810444
$ws.Range("J16").Value = 1.073001755382389
$ws.Range("K16").Value = 1.065084034084724
$ws.Range("L16").Value = 1.083788381962338
$ws.Range("M16").Value = 1.090489256542724
$ws.Range("N16").Value = 1.074525540950053

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.06746262041643
$ws.Range("D17").Value = 1.06216828925276
$ws.Range("E17").Value = 1.08108457579255
$ws.Range("F17").Value = 1.087862298948141
$ws.Range("I17").Value = 1.051171550594682
$ws.Range("J17").Value = 1.073440145989355
$ws.Range("K17").Value = 1.065445399342589
$ws.Range("L17").Value = 1.084300474666388
$ws.Range("M17").Value = 1.091056832175506
$ws.Range("N17").Value = 1.074964554121965

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.067796212873001
$ws.Range("D18").Value = 1.062420388574467
$ws.Range("E18").Value = 1.081424066444807
$ws.Range("F18").Value = 1.088233984946521
$ws.Range("I18").Value = 1.051278975851378
$ws.Range("J18").Value = 1.073695620561082
$ws.Range("K18").Value = 1.065655942911719
$ws.Range("L18").Value = 1.084598972038928
$ws.Range("M18").Value = 1.091387720352803
$ws.Range("N18").Value = 1.075220391496888

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.067909928686153
$ws.Range("D19").Value = 1.062506320011052
$ws.Range("E19").Value = 1.081539800801109
$ws.Range("F19").Value = 1.088360701386784
$ws.Range("I19").Value = 1.051315564975605
$ws.Range("J19").Value = 1.073782691753564
$ws.Range("K19").Value = 1.06572769317695
$ws.Range("L19").Value = 1.08470071862401
$ws.Range("M19").Value = 1.091500516230447
$ws.Range("N19").Value = 1.07530758634045

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.06740124401369
$ws.Range("D20").Value = 1.0621219041903
$ws.Range("E20").Value = 1.081022117981254
$ws.Range("F20").Value = 1.087793921081173
$ws.Range("I20").Value = 1.051151771466959
$ws.Range("J20").Value = 1.073393134791413
$ws.Range("K20").Value = 1.06540665258059
$ws.Range("L20").Value = 1.0842455524607
$ws.Range("M20").Value = 1.090995954228864
$ws.Range("N20").Value = 1.07491747616273

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.065745305216132
$ws.Range("D21").Value = 1.060870173570895
$ws.Range("E21").Value = 1.079337440481843
$ws.Range("F21").Value = 1.085949917236397
$ws.Range("I21").Value = 1.050616470672249
$ws.Range("J21").Value = 1.072123915242665
$ws.Range("K21").Value = 1.064360148896724
$ws.Range("L21").Value = 1.082763424061033
$ws.Range("M21").Value = 1.089353569074239
$ws.Range("N21").Value = 1.073646454176642

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.064702455722831
$ws.Range("D22").Value = 1.060081632728405
$ws.Range("E22").Value = 1.078276918622072
$ws.Range("F22").Value = 1.084789441624668
$ws.Range("I22").Value = 1.050277750302088
$ws.Range("J22").Value = 1.071323779651954
$ws.Range("K22").Value = 1.063700021142738
$ws.Range("L22").Value = 1.081829720713203
$ws.Range("M22").Value = 1.088319358030123
$ws.Range("N22").Value = 1.072845182301615

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.065255453409852
$ws.Range("D23").Value = 1.06049979978831
$ws.Range("E23").Value = 1.078839247027931
$ws.Range("F23").Value = 1.08540473674896
$ws.Range("I23").Value = 1.050457517685093
$ws.Range("J23").Value = 1.071748150389983
$ws.Range("K23").Value = 1.064050172801546
$ws.Range("L23").Value = 1.082324870451369
$ws.Range("M23").Value = 1.088867765154355
$ws.Range("N23").Value = 1.073270155694768

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.06742897793029
$ws.Range("D24").Value = 1.062142864115178
$ws.Range("E24").Value = 1.08105034040465
$ws.Range("F24").Value = 1.087824818445003
$ws.Range("I24").Value = 1.05116070953586
$ws.Range("J24").Value = 1.073414377842978
$ws.Range("K24").Value = 1.065424161299839
$ws.Range("L24").Value = 1.084270370052299
$ws.Range("M24").Value = 1.091023462875231
$ws.Range("N24").Value = 1.074938749381865

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.06994412784079
$ws.Range("D25").Value = 1.064043073625221
$ws.Range("E25").Value = 1.08361080211424
$ws.Range("F25").Value = 1.090628791553356
$ws.Range("I25").Value = 1.051967437353373
$ws.Range("J25").Value = 1.075338892383888
$ws.Range("K25").Value = 1.067009410434401
$ws.Range("L25").Value = 1.086520286216607
$ws.Range("M25").Value = 1.093518424860632
$ws.Range("N25").Value = 1.076865996954169
